$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "(`"Ass Whuppin'`", ['{1}{W}{B}', 'Sorcery', 'Destroy target silver-bordered permanent in any game you can see from your seat.'])"

$ws.Range("A3:A5").EntireRow.Delete()
